$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)
$shape.TextFrame.TextRange.Text = "Removes bottom half of the stomach`rDoes not remove all lymph nodes"
